# Auto-generated Excel COM edit script
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # 展览
$ws3 = $wb.Worksheets.Item(3)  # 本地生活
$ws4 = $wb.Worksheets.Item(4)  # 全部类型

# ---- 展览 ----
$ws1.Range("F5").Value = 308
$ws1.Range("F8").Value = 2001
$ws1.Range("F11").Value = 32
$ws1.Range("F12").Value = 1594
$ws1.Range("F13").Value = 1594
$ws1.Range("F14").Value = 1323
$ws1.Range("F17").Value = 179
$ws1.Range("F20").Value = 448
$ws1.Range("B22").Value = "'2024-09-22"
$ws1.Range("C22").Value = "北京·地狱双ip同人ONLY展"
$ws1.Range("D22").Value = "双桥中路50号院 E50艺术园区"
$ws1.Range("E22").Value = "2024.09.22 10:30-09.22 16:00"
$ws1.Range("F22").Value = 142
$ws1.Range("G22").Value = 105
$ws1.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=90931"
$ws1.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202408/c6ObwO4C1724055713128.jpeg"
$ws1.Range("B23").Value = "'2024-10-01"
$ws1.Range("C23").Value = "北京·IDO动漫游戏嘉年华47th"
$ws1.Range("D23").Value = "亦庄荣昌东街6号 北京亦创国际会展中心"
$ws1.Range("E23").Value = "2024.10.01 09:30-10.02 17:00"
$ws1.Range("F23").Value = 7032
$ws1.Range("G23").Value = 5
$ws1.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=83826"
$ws1.Range("I23").Value = "//i0.hdslb.com/bfs/openplatform/202405/JL6boAFV1716882961702.jpeg"
$ws1.Range("F24").Value = 7032
$ws1.Range("C25").Value = "北京·第19届IJOY漫展xCGF游戏节"
$ws1.Range("D25").Value = "天辰东路7号 北京国家会议中心"
$ws1.Range("E25").Value = "2024.10.01 09:00-10.02 17:00"
$ws1.Range("F25").Value = 7605
$ws1.Range("G25").Value = 8.800000000000001
$ws1.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=84127"
$ws1.Range("I25").Value = "//i0.hdslb.com/bfs/openplatform/202405/iR6rV5311717039317028.jpeg"
$ws1.Range("C26").Value = "北京·第19届IJOY漫展【Pile专场见面会】"
$ws1.Range("E26").Value = "2024.10.01 14:50-10.01 16:30"
$ws1.Range("F26").Value = 36
$ws1.Range("G26").Value = 458
$ws1.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=91560"
$ws1.Range("I26").Value = "//i1.hdslb.com/bfs/openplatform/202408/mBtVCKBp1724927832154.jpeg"
$ws1.Range("C27").Value = "北京·第19届IJOY漫展【广播剧《伪装者》专场见面会】"
$ws1.Range("E27").Value = "2024.10.01 11:00-10.01 15:30"
$ws1.Range("F27").Value = 2
$ws1.Range("G27").Value = 388
$ws1.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=91771"
$ws1.Range("I27").Value = "//i1.hdslb.com/bfs/openplatform/202409/6yy1mOUn1725334481066.jpeg"
$ws1.Range("C28").Value = "北京·第五人格同人only同人3.0"
$ws1.Range("D28").Value = "永外高庄138号 北京大红门国际会展中心"
$ws1.Range("E28").Value = "2024.10.01 10:00-10.01 17:00"
$ws1.Range("F28").Value = 178
$ws1.Range("G28").Value = 60
$ws1.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=90653"
$ws1.Range("I28").Value = "//i2.hdslb.com/bfs/openplatform/202408/UsBZWtUX1723532208881.jpeg"
$ws1.Range("C29").Value = "北京·配音演员 金弦 专场活动"
$ws1.Range("D29").Value = "亦庄荣昌东街6号 北京亦创国际会展中心"
$ws1.Range("E29").Value = "2024.10.01 10:30-10.01 13:30"
$ws1.Range("F29").Value = 491
$ws1.Range("G29").Value = "已售罄"
$ws1.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=91068"
$ws1.Range("I29").Value = "//i1.hdslb.com/bfs/openplatform/202408/vJRCM3vg1724226523747.jpeg"
$ws1.Range("B30").Value = "'2024-10-02"
$ws1.Range("C30").Value = "北京·人气声优 内田秀 专场活动"
$ws1.Range("E30").Value = "2024.10.02 13:55-10.02 17:10"
$ws1.Range("F30").Value = 77
$ws1.Range("G30").Value = 458
$ws1.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=91678"
$ws1.Range("I30").Value = "//i0.hdslb.com/bfs/openplatform/202409/0aUkHD511725260741169.png"
$ws1.Range("C31").Value = "北京·人气声优 小林爱香 专场活动"
$ws1.Range("E31").Value = "2024.10.02 12:50-10.02 16:40"
$ws1.Range("F31").Value = 216
$ws1.Range("G31").Value = "已售罄"
$ws1.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=91117"
$ws1.Range("I31").Value = "//i2.hdslb.com/bfs/openplatform/202408/nuqS5Gd11724309352207.png"
$ws1.Range("C32").Value = "北京·人气声优 青山渚 专场活动"
$ws1.Range("E32").Value = "2024.10.02 11:50-10.02 15:40"
$ws1.Range("F32").Value = 247
$ws1.Range("G32").Value = 458
$ws1.Range("H32").Value = "https://show.bilibili.com/platform/detail.html?id=91249"
$ws1.Range("I32").Value = "//i2.hdslb.com/bfs/openplatform/202408/xHqpdFa41724641733192.png"
$ws1.Range("C33").Value = "北京·广播剧《西东》专场活动"
$ws1.Range("E33").Value = "2024.10.02 10:30-10.02 13:50"
$ws1.Range("F33").Value = 50
$ws1.Range("G33").Value = 300
$ws1.Range("H33").Value = "https://show.bilibili.com/platform/detail.html?id=91844"
$ws1.Range("I33").Value = "//i0.hdslb.com/bfs/openplatform/202409/rhW19Bur1725418971000.png"
$ws1.Range("F38").Value = 1380
$ws1.Range("F41").Value = 278
$ws1.Range("F42").Value = 684
$ws1.Range("F43").Value = 21

# ---- 本地生活 ----
$ws3.Range("F3").Value = 2565

# ---- 全部类型 ----
$ws4.Range("F9").Value = 308
$ws4.Range("F12").Value = 2001
$ws4.Range("F14").Value = 32
$ws4.Range("F15").Value = 1594
$ws4.Range("F16").Value = 1594
$ws4.Range("F17").Value = 1323
$ws4.Range("F20").Value = 448
$ws4.Range("F24").Value = 7032
$ws4.Range("F25").Value = 7032
$ws4.Range("F26").Value = 7605
$ws4.Range("F30").Value = 247
$ws4.Range("F34").Value = 1380
$ws4.Range("F38").Value = 278
$ws4.Range("F41").Value = 684
